# Update "Unidades Pedido" (L) and "Diferencia Stock" (M) columns for the
# affected item rows, then refresh the total rows (C88 = Total_Unidades,
# C99 = Total_Ajuste_Stock) to reflect the new sums.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_14")

$changes = @{
    5  = @(2, -1)
    10 = @(2, 1)
    15 = @(5, 1)
    23 = @(12, 2)
    25 = @(5, 1)
    26 = @(2, 1)
    29 = @(3, 1)
    30 = @(8, 2)
    34 = @(2, 1)
    35 = @(5, 1)
    36 = @(2, 1)
    40 = @(2, 1)
    47 = @(4, 1)
    48 = @(2, 1)
    70 = @(2, 1)
    74 = @(8, 2)
    81 = @(11, 2)
    82 = @(7, 1)
    85 = @(9, 2)
}

foreach ($row in $changes.Keys) {
    $vals = $changes[$row]
    $ws.Range("L$row").Value = $vals[0]
    $ws.Range("M$row").Value = $vals[1]
}

$ws.Range("C88").Value = 236
$ws.Range("C99").Value = 22
